$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "53"
$ws.Range("C3").NumberFormat = "General"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "5"
$ws.Range("E3").NumberFormat = "General"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "2"
$ws.Range("F3").NumberFormat = "General"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "17"
$ws.Range("C4").NumberFormat = "General"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "9"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1"
$ws.Range("E4").NumberFormat = "General"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "1"
$ws.Range("F4").NumberFormat = "General"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "15"
$ws.Range("C5").NumberFormat = "General"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "11"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "1"
$ws.Range("F5").NumberFormat = "General"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "5"
$ws.Range("C6").NumberFormat = "General"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0"
$ws.Range("E6").NumberFormat = "General"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "0"
$ws.Range("C7").NumberFormat = "General"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "0"
$ws.Range("F7").NumberFormat = "General"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "1"
$ws.Range("C8").NumberFormat = "General"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0"
$ws.Range("E8").NumberFormat = "General"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "0"
$ws.Range("F8").NumberFormat = "General"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "33"
$ws.Range("C9").NumberFormat = "General"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "12"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1"
$ws.Range("E9").NumberFormat = "General"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "4"
$ws.Range("F9").NumberFormat = "General"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "4"
$ws.Range("C10").NumberFormat = "General"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0"
$ws.Range("E10").NumberFormat = "General"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "0"
$ws.Range("F10").NumberFormat = "General"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "5"
$ws.Range("C11").NumberFormat = "General"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "0"
$ws.Range("F11").NumberFormat = "General"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "12"
$ws.Range("C12").NumberFormat = "General"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "10"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1"
$ws.Range("E12").NumberFormat = "General"

Write-Host "Updated Pat Cummins batting activity rows 3-12"